$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 - "What is an ML Pipeline?"
#   "Testing ideas and hypothesis easily and quickly"
#     -> "Testing ideas and hypotheses easily and quickly"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4body = $s4.Shapes.Item(1).TextFrame.TextRange
$s4count = $s4body.Paragraphs().Count
for ($i = 1; $i -le $s4count; $i++) {
    $para = $s4body.Paragraphs($i, 1)
    if ($para.Text -eq "Testing ideas and hypothesis easily and quickly") {
        # Two-stage replace so the whole paragraph collapses back into a
        # single run (a direct in-place edit would otherwise produce a
        # common-prefix/suffix split into multiple runs).
        $para.Text = "zzz_tmp_zzz"
        $para.Text = "Testing ideas and hypotheses easily and quickly"
        break
    }
}

# ---------------------------------------------------------------------------
# Slide 12 - "Feature Creation"
#   insert "Data comes with fields or columns (if it's even structured), not
#     features" before "Common Features"
#   insert "How are you handling imputation of missing values?" after
#     "Spatial"
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12body = $s12.Shapes.Item(1).TextFrame.TextRange

$firstPara = $s12body.Paragraphs(1, 1)
$apos = [char]0x2019
$firstPara.InsertBefore("Data comes with fields or columns (if it${apos}s even structured), not features`r")

$s12count = $s12body.Paragraphs().Count
$lastPara = $s12body.Paragraphs($s12count, 1)
$lastPara.InsertAfter("`rHow are you handling imputation of missing values?")
$s12count2 = $s12body.Paragraphs().Count
$newPara = $s12body.Paragraphs($s12count2, 1)
$newPara.IndentLevel = 1

# ---------------------------------------------------------------------------
# Slide 15 - "Deployment"
#   insert "Model monitoring" before "Re-training"
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$s15body = $s15.Shapes.Item(1).TextFrame.TextRange
$s15firstPara = $s15body.Paragraphs(1, 1)
$s15firstPara.InsertBefore("Model monitoring`r")

# ---------------------------------------------------------------------------
# Slide 17 - "Best Practices"
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$s17body = $s17.Shapes.Item(1).TextFrame.TextRange
$apos17 = [char]0x2019

$configPara = $s17body.Paragraphs(1, 1)

# Merge ", " + "json" + ", " (chars 19-26) into a single ", json, " run
# (taking the formatting of the leading ", " run, which has no err flag).
$firstComma = $configPara.Characters(19, 2)
$firstComma.InsertAfter("json, ")
$oldJsonComma = $configPara.Characters(27, 6)
$oldJsonComma.Text = ""

# Merge "Config" + " files (" (chars 1-14) into a single "Config files ("
# run (taking the formatting of the " files (" run, which has no err flag).
$configRun = $configPara.Characters(1, 6)
$configRun.Text = ""
$filesRun = $configPara.Characters(1, 8)
$filesRun.InsertBefore("Config")

# Insert the two new top paragraphs before the "Config files (...)" line.
$configPara.InsertBefore("Draw a diagram of the pipeline: `rWhat function runs each step? What are the inputs? What are the outputs?`r")

# "What function runs..." paragraph is now #2 -> demote to level 1 (lvl="1").
$whatFuncPara = $s17body.Paragraphs(2, 1)
$whatFuncPara.IndentLevel = 2

# "Config files (...)" paragraph is now #3; insert the 3 new paragraphs
# right after it and before "Store models as pickles".
$configFilesPara = $s17body.Paragraphs(3, 1)
$configFilesPara.InsertAfter("`rMake each step modular and extensible so it can easily be re-used`rBuild a simple, end-to-end version first, then add more functionality`rThink about how you${apos17}ll store outputs:")

# The trailing 4 paragraphs (Store models / Store predictions / Store
# evaluation / Sample results schema) are now #7-#10; demote them to
# level 1 (lvl="1").
$s17count = $s17body.Paragraphs().Count
for ($i = $s17count - 3; $i -le $s17count; $i++) {
    $pp = $s17body.Paragraphs($i, 1)
    $pp.IndentLevel = 2
}
